$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a "codeforiati:group-name" column (D) and a
# "codeforiati:group-code" column (E) whose order needs to be swapped
# (both the header text and every data row's values), so that the
# group-code column comes before the group-name column.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$colD = $ws.Range("D1:D$lastRow")
$colE = $ws.Range("E1:E$lastRow")

$dValues = $colD.Value2
$eValues = $colE.Value2

$colD.Value2 = $eValues
$colE.Value2 = $dValues
